$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "0001_slr0611_right"
$ws.Range("C4").Value = "this is a row with an empty attachment file"
$ws.Range("D4").Value = "this is a row with an empty attachment file"

$ws.Range("B5").Value = "NC_014139.gbk"
$ws.Range("D5").Value = "This is a row with an empty description"
$ws.Range("A5").Value = "0003_slr0613_left"

$ws.Range("A5").Select()

$ws.PageSetup.Orientation = 1
